$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: set to old row 39's data
$ws.Range("F38").Value = "Marines"
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = "AS Kigali"
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 2.99
$ws.Range("K38").Value = "10/10/2023 02:12"
$ws.Range("L38").Value = 2.75
$ws.Range("M38").Value = "11/10/2023 14:57"
$ws.Range("N38").Value = 2.78
$ws.Range("O38").Value = "10/10/2023 02:12"
$ws.Range("P38").Value = 2.88
$ws.Range("Q38").Value = "11/10/2023 14:57"
$ws.Range("R38").Value = 2.21
$ws.Range("S38").Value = "10/10/2023 02:12"
$ws.Range("T38").Value = 2.54
$ws.Range("U38").Value = "11/10/2023 14:57"
$ws.Range("V38").Value = "https://www.betexplorer.com/football/rwanda/premier-league/marines-as-kigali/QHn271so/"

# Row 39: set to old row 38's data
$ws.Range("F39").Value = "Amagaju"
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = "Kiyovu"
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 2.81
$ws.Range("K39").Value = "11/10/2023 11:11"
$ws.Range("L39").Value = 3.14
$ws.Range("M39").Value = "11/10/2023 14:10"
$ws.Range("N39").Value = 2.82
$ws.Range("O39").Value = "11/10/2023 11:11"
$ws.Range("P39").Value = 2.93
$ws.Range("Q39").Value = "11/10/2023 14:10"
$ws.Range("R39").Value = 2.41
$ws.Range("S39").Value = "11/10/2023 11:11"
$ws.Range("T39").Value = 2.25
$ws.Range("U39").Value = "11/10/2023 14:10"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/rwanda/premier-league/amagaju-kiyovu/lMQVEqkB/"

# Row 46: set to old row 48's data
$ws.Range("F46").Value = "Etoile de L'Est"
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = "Muhazi United"
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 2.45
$ws.Range("K46").Value = "15/10/2023 13:12"
$ws.Range("L46").Value = 2.42
$ws.Range("M46").Value = "15/10/2023 14:45"
$ws.Range("N46").Value = 2.77
$ws.Range("O46").Value = "15/10/2023 13:12"
$ws.Range("P46").Value = 2.8
$ws.Range("Q46").Value = "15/10/2023 14:45"
$ws.Range("R46").Value = 2.91
$ws.Range("S46").Value = "15/10/2023 13:12"
$ws.Range("T46").Value = 2.98
$ws.Range("U46").Value = "15/10/2023 14:45"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/rwanda/premier-league/etoile-de-l-est-muhazi-united/tjrSEb71/"

# Row 47: set to old row 46's data
$ws.Range("F47").Value = "Musanze"
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = "Rayon Sport"
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 3.66
$ws.Range("K47").Value = "14/10/2023 02:12"
$ws.Range("L47").Value = 3.04
$ws.Range("M47").Value = "15/10/2023 14:58"
$ws.Range("N47").Value = 2.93
$ws.Range("O47").Value = "14/10/2023 02:12"
$ws.Range("P47").Value = 2.4
$ws.Range("Q47").Value = "15/10/2023 14:58"
$ws.Range("R47").Value = 1.87
$ws.Range("S47").Value = "14/10/2023 02:12"
$ws.Range("T47").Value = 2.79
$ws.Range("U47").Value = "15/10/2023 14:58"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/rwanda/premier-league/musanze-rayon-sport/bRhNFvhe/"

# Row 48: set to old row 47's data
$ws.Range("F48").Value = "AS Kigali"
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = "Police"
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 1.97
$ws.Range("K48").Value = "14/10/2023 02:12"
$ws.Range("L48").Value = 2.07
$ws.Range("M48").Value = "15/10/2023 11:02"
$ws.Range("N48").Value = 2.82
$ws.Range("O48").Value = "14/10/2023 02:12"
$ws.Range("P48").Value = 2.81
$ws.Range("Q48").Value = "15/10/2023 13:01"
$ws.Range("R48").Value = 3.49
$ws.Range("S48").Value = "14/10/2023 02:12"
$ws.Range("T48").Value = 3.7
$ws.Range("U48").Value = "15/10/2023 11:02"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/rwanda/premier-league/as-kigali-police/KzfJGKxk/"

# Row 52: set to old row 54's data
$ws.Range("F52").Value = "Amagaju"
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = "Gorilla"
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 1.97
$ws.Range("K52").Value = "21/10/2023 02:13"
$ws.Range("L52").Value = 2.09
$ws.Range("M52").Value = "22/10/2023 14:25"
$ws.Range("N52").Value = 2.94
$ws.Range("O52").Value = "21/10/2023 02:13"
$ws.Range("P52").Value = 2.84
$ws.Range("Q52").Value = "22/10/2023 14:25"
$ws.Range("R52").Value = 3.32
$ws.Range("S52").Value = "21/10/2023 02:13"
$ws.Range("T52").Value = 3.63
$ws.Range("U52").Value = "22/10/2023 14:25"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/rwanda/premier-league/amagaju-gorilla/CUIdUw7l/"

# Row 53: set to old row 52's data
$ws.Range("F53").Value = "Etincelles"
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = "APR"
$ws.Range("I53").Value = 3
$ws.Range("J53").Value = 4.19
$ws.Range("K53").Value = "21/10/2023 07:28"
$ws.Range("L53").Value = 4.19
$ws.Range("M53").Value = "21/10/2023 07:28"
$ws.Range("N53").Value = 3.58
$ws.Range("O53").Value = "21/10/2023 13:02"
$ws.Range("P53").Value = 3.58
$ws.Range("Q53").Value = "21/10/2023 13:02"
$ws.Range("R53").Value = 1.68
$ws.Range("S53").Value = "21/10/2023 07:28"
$ws.Range("T53").Value = 1.68
$ws.Range("U53").Value = "21/10/2023 07:28"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/rwanda/premier-league/etincelles-apr/rFSQZaqR/"

# Row 54: set to old row 53's data
$ws.Range("F54").Value = "Etoile de L'Est"
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = "AS Kigali"
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 3
$ws.Range("K54").Value = "21/10/2023 14:10"
$ws.Range("L54").Value = 3
$ws.Range("M54").Value = "21/10/2023 14:10"
$ws.Range("N54").Value = 3.15
$ws.Range("O54").Value = "21/10/2023 14:10"
$ws.Range("P54").Value = 3.15
$ws.Range("Q54").Value = "21/10/2023 14:10"
$ws.Range("R54").Value = 2.2
$ws.Range("S54").Value = "21/10/2023 14:10"
$ws.Range("T54").Value = 2.2
$ws.Range("U54").Value = "21/10/2023 14:10"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/rwanda/premier-league/etoile-de-l-est-as-kigali/Yq75SHy1/"

# Row 68: set to old row 71's data
$ws.Range("F68").Value = "Etincelles"
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = "Rayon Sport"
$ws.Range("I68").Value = 1
$ws.Range("J68").Value = 4.18
$ws.Range("K68").Value = "11/11/2023 03:13"
$ws.Range("L68").Value = 3.12
$ws.Range("M68").Value = "25/11/2023 13:56"
$ws.Range("N68").Value = 3.15
$ws.Range("O68").Value = "11/11/2023 03:13"
$ws.Range("P68").Value = 2.74
$ws.Range("Q68").Value = "25/11/2023 13:56"
$ws.Range("R68").Value = 1.77
$ws.Range("S68").Value = "11/11/2023 03:13"
$ws.Range("T68").Value = 2.38
$ws.Range("U68").Value = "25/11/2023 13:56"
$ws.Range("V68").Value = "https://www.betexplorer.com/football/rwanda/premier-league/etincelles-rayon-sport/tpRMiY4C/"

# Row 70: set to old row 68's data
$ws.Range("F70").Value = "APR"
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = "AS Kigali"
$ws.Range("I70").Value = 1
$ws.Range("J70").Value = 1.7
$ws.Range("K70").Value = "11/11/2023 03:13"
$ws.Range("L70").Value = 1.77
$ws.Range("M70").Value = "25/11/2023 12:54"
$ws.Range("N70").Value = 3.13
$ws.Range("O70").Value = "11/11/2023 03:13"
$ws.Range("P70").Value = 3.04
$ws.Range("Q70").Value = "25/11/2023 12:54"
$ws.Range("R70").Value = 4.63
$ws.Range("S70").Value = "11/11/2023 03:13"
$ws.Range("T70").Value = 4.74
$ws.Range("U70").Value = "25/11/2023 12:54"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/rwanda/premier-league/apr-as-kigali/nqgfZPlo/"

# Row 71: set to old row 70's data
$ws.Range("F71").Value = "Bugesera"
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = "Marines"
$ws.Range("I71").Value = 1
$ws.Range("J71").Value = 1.99
$ws.Range("K71").Value = "11/11/2023 03:13"
$ws.Range("L71").Value = 2.4
$ws.Range("M71").Value = "25/11/2023 13:54"
$ws.Range("N71").Value = 3.07
$ws.Range("O71").Value = "11/11/2023 03:13"
$ws.Range("P71").Value = 2.91
$ws.Range("Q71").Value = "25/11/2023 13:54"
$ws.Range("R71").Value = 3.42
$ws.Range("S71").Value = "11/11/2023 03:13"
$ws.Range("T71").Value = 2.9
$ws.Range("U71").Value = "25/11/2023 13:54"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/rwanda/premier-league/bugesera-marines/2gSIhEk6/"

# New row 74: append new match data (shift dimension to A1:V74)
$ws.Range("A73:V73").Copy()
$ws.Range("A74:V74").PasteSpecial(-4122)
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "rwanda"
$ws.Range("C74").Value = "premier-league"
$ws.Range("D74").Value = "2023-2024"
$ws.Range("E74").Value = 45258.58333333334
$ws.Range("F74").Value = "Police"
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = "Rayon Sport"
$ws.Range("I74").Value = 2
$ws.Range("J74").Value = 2.73
$ws.Range("K74").Value = "28/11/2023 03:12"
$ws.Range("L74").Value = 2.6
$ws.Range("M74").Value = "28/11/2023 13:45"
$ws.Range("N74").Value = 2.76
$ws.Range("O74").Value = "28/11/2023 03:12"
$ws.Range("P74").Value = 2.61
$ws.Range("Q74").Value = "28/11/2023 13:45"
$ws.Range("R74").Value = 2.62
$ws.Range("S74").Value = "28/11/2023 03:12"
$ws.Range("T74").Value = 2.97
$ws.Range("U74").Value = "28/11/2023 13:45"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/rwanda/premier-league/police-rayon-sport/KMAts8XA/"
